$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.247.62'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').Value = '2.273.09'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '299.86'
$ws.Range('E5').Value = '  -1.42%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '96.33'
$ws.Range('E6').Value = '  -3.09%  '
$ws.Range('E7').Value = '  -1.75%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -1.84%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '33.33'
$ws.Range('E10').Value = '  -3.19%  '
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '48.30'
$ws.Range('E12').Value = '  -6.75%  '
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.67'
$ws.Range('E14').Value = '  -1.14%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.59'
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('D16').Value = '2.627.69'
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('D17').Value = '2.290.03'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.785'
$ws.Range('E18').Value = '  -4.32%  '
$ws.Range('D19').Value = '42.190.42'
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.73'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('D21').Value = '0.0₃0891'
$ws.Range('E21').Value = '  -1.22%  '
$ws.Range('E22').Value = '  -1.45%  '
$ws.Range('E23').Value = '  -3.83%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '235.05'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.98'
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.46'
$ws.Range('E27').Value = '  -2.84%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '24.02'
$ws.Range('E28').Value = '  -4.81%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '168.98'
$ws.Range('E29').Value = '  +4.07%  '
$ws.Range('E30').Value = '  -4.91%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '9.21'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  -2.65%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  -2.68%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.55'
$ws.Range('E35').Value = '  -1.77%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '16.63'
$ws.Range('E36').Value = '  -2.28%  '
$ws.Range('E37').Value = '  -4.74%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0688'
$ws.Range('E38').Value = '  -3.95%  '
$ws.Range('E39').Value = '  -3.44%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0990'
$ws.Range('E40').Value = '  -1.53%  '
$ws.Range('E41').Value = '  -2.42%  '
$ws.Range('E42').Value = '  -4.58%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.41'
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('D44').Value = '1.973.22'
$ws.Range('E44').Value = '  -1.01%  '
$ws.Range('E45').Value = '  -0.82%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '17.52'
$ws.Range('E46').Value = '  -6.50%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.57'
$ws.Range('E47').Value = '  -6.72%  '
$ws.Range('E48').Value = '  -2.95%  '
$ws.Range('D49').Value = '2.498.31'
$ws.Range('E49').Value = '  -1.53%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '52.37'
$ws.Range('E50').Value = '  -5.57%  '
$ws.Range('E51').Value = '  -1.00%  '
